$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.203.80'
$ws.Range("E2").Value = '  +0.00%  '
$ws.Range("D3").Value = '2.490.47'
$ws.Range("E3").Value = '  +0.08%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.32'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.30%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.31'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.47%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D9").Value = '2.489.77'
$ws.Range("E9").Value = '  +0.09%  '
$ws.Range("E10").Value = '  +0.57%  '
$ws.Range("E11").Value = '  +0.13%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.93'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.39%  '
$ws.Range("E13").Value = '  -2.04%  '
$ws.Range("E14").Value = '  +0.37%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '25.41'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.30%  '
$ws.Range("D16").Value = '67.270.13'
$ws.Range("E16").Value = '  +0.73%  '
$ws.Range("E17").Value = '  -1.99%  '
$ws.Range("D18").Value = '2.482.41'
$ws.Range("E18").Value = '  +1.19%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.00'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -6.32%  '
$ws.Range("E20").Value = '  -4.86%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '349.25'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.19%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.02'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.32%  '
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '68.62'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.13%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.22'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.75%  '
$ws.Range("E26").Value = '  -3.15%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.28'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.03%  '
$ws.Range("E28").Value = '  +0.17%  '
$ws.Range("E29").Value = '  +0.11%  '
$ws.Range("E30").Value = '  -3.66%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '509.65'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.60%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.77'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.05%  '
$ws.Range("E33").Value = '  -3.19%  '
$ws.Range("E34").Value = '  -3.97%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.08%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '159.87'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.86%  '
$ws.Range("E37").Value = '  -7.55%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.71'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.76%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.25'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.70%  '
$ws.Range("E40").Value = '  -6.02%  '
$ws.Range("E41").Value = '  -2.59%  '
$ws.Range("E42").Value = '  -0.05%  '
$ws.Range("B43").Value = 'PolygonEcosystemToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.328'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.20%  '
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.82'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.94%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.37'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.23%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '38.80'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.53%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '142.51'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.30%  '
$ws.Range("E48").Value = '  -4.73%  '
$ws.Range("E49").Value = '  -4.49%  '
$ws.Range("E50").Value = '  -6.35%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0732'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.77%  '
